$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamadmin', 'admin: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamadmin', 'maintain: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamadmin', 'push: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamadmin', 'triage: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamadmin', 'pull: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamwrite', 'admin: False'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamwrite', 'maintain: False'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamwrite', 'push: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamwrite', 'triage: True'),
    @('github-gk-aks/firstgithubrepo', 'Team: firstgithubrepo-teamwrite', 'pull: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamadmin', 'admin: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamadmin', 'maintain: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamadmin', 'push: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamadmin', 'triage: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamadmin', 'pull: True'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamread', 'admin: False'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamread', 'maintain: False'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamread', 'push: False'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamread', 'triage: False'),
    @('github-gk-aks/secondgithubrepo', 'Team: secondgithubrepo-teamread', 'pull: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamadmin', 'admin: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamadmin', 'maintain: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamadmin', 'push: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamadmin', 'triage: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamadmin', 'pull: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamread', 'admin: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamread', 'maintain: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamread', 'push: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamread', 'triage: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamread', 'pull: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamwrite', 'admin: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamwrite', 'maintain: False'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamwrite', 'push: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamwrite', 'triage: True'),
    @('github-gk-aks/thirdgithubrepo', 'Team: thirdgithubrepo-teamwrite', 'pull: True'),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
